$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("February")

$ws.Cells.Item(55,2).Value = 288
$ws.Cells.Item(55,3).Value = 221
$ws.Cells.Item(55,4).Value = 67
$ws.Cells.Item(55,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(55,7).Value = "1.30 : 1"

$ws.Cells.Item(54,2).Value = 16
$ws.Cells.Item(54,3).Value = 225
$ws.Cells.Item(54,4).Value = -209
$ws.Cells.Item(54,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(54,7).Value = "0.07 : 1"

$ws.Cells.Item(53,2).Value = 138
$ws.Cells.Item(53,3).Value = 204
$ws.Cells.Item(53,4).Value = -66
$ws.Cells.Item(53,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(53,7).Value = "0.68 : 1"

$ws.Cells.Item(52,2).Value = 367
$ws.Cells.Item(52,3).Value = 345
$ws.Cells.Item(52,4).Value = 22
$ws.Cells.Item(52,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(52,7).Value = "1.06 : 1"

$ws.Cells.Item(51,2).Value = 204
$ws.Cells.Item(51,3).Value = 166
$ws.Cells.Item(51,4).Value = 38
$ws.Cells.Item(51,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(51,7).Value = "1.23 : 1"

$ws.Cells.Item(50,2).Value = 897
$ws.Cells.Item(50,3).Value = 531
$ws.Cells.Item(50,4).Value = 366
$ws.Cells.Item(50,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(50,7).Value = "1.69 : 1"

$ws.Cells.Item(49,2).Value = 539
$ws.Cells.Item(49,3).Value = 214
$ws.Cells.Item(49,4).Value = 325
$ws.Cells.Item(49,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(49,7).Value = "2.52 : 1"

$ws.Cells.Item(48,2).Value = 263
$ws.Cells.Item(48,3).Value = 638
$ws.Cells.Item(48,4).Value = -375
$ws.Cells.Item(48,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(48,7).Value = "0.41 : 1"

$ws.Cells.Item(47,2).Value = 1143
$ws.Cells.Item(47,3).Value = 605
$ws.Cells.Item(47,4).Value = 538
$ws.Cells.Item(47,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(47,7).Value = "1.89 : 1"

$ws.Cells.Item(46,2).Value = 608
$ws.Cells.Item(46,3).Value = 582
$ws.Cells.Item(46,4).Value = 26
$ws.Cells.Item(46,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(46,7).Value = "1.04 : 1"

$ws.Cells.Item(45,2).Value = 66
$ws.Cells.Item(45,3).Value = 183
$ws.Cells.Item(45,4).Value = -117
$ws.Cells.Item(45,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(45,7).Value = "0.36 : 1"

$ws.Cells.Item(44,2).Value = 77
$ws.Cells.Item(44,3).Value = 83
$ws.Cells.Item(44,4).Value = -6
$ws.Cells.Item(44,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(44,7).Value = "0.93 : 1"

$ws.Cells.Item(43,2).Value = 0
$ws.Cells.Item(43,3).Value = 0
$ws.Cells.Item(43,4).Value = 0

$ws.Cells.Item(42,2).Value = 12
$ws.Cells.Item(42,3).Value = 37
$ws.Cells.Item(42,4).Value = -25
$ws.Cells.Item(42,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(42,7).Value = "0.32 : 1"

$ws.Cells.Item(41,2).Value = 3
$ws.Cells.Item(41,3).Value = 36
$ws.Cells.Item(41,4).Value = -33
$ws.Cells.Item(41,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(41,7).Value = "0.08 : 1"

$ws.Cells.Item(40,2).Value = 96
$ws.Cells.Item(40,3).Value = 92
$ws.Cells.Item(40,4).Value = 4
$ws.Cells.Item(40,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(40,7).Value = "1.04 : 1"

$ws.Cells.Item(39,2).Value = 21
$ws.Cells.Item(39,3).Value = 114
$ws.Cells.Item(39,4).Value = -93
$ws.Cells.Item(39,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(39,7).Value = "0.18 : 1"

$ws.Cells.Item(38,2).Value = 24
$ws.Cells.Item(38,3).Value = 188
$ws.Cells.Item(38,4).Value = -164
$ws.Cells.Item(38,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(38,7).Value = "0.13 : 1"

$ws.Cells.Item(37,2).Value = 440
$ws.Cells.Item(37,3).Value = 306
$ws.Cells.Item(37,4).Value = 134
$ws.Cells.Item(37,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(37,7).Value = "1.44 : 1"

$ws.Cells.Item(36,2).Value = 235
$ws.Cells.Item(36,3).Value = 490
$ws.Cells.Item(36,4).Value = -255
$ws.Cells.Item(36,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(36,7).Value = "0.48 : 1"

$ws.Cells.Item(35,2).Value = 832
$ws.Cells.Item(35,3).Value = 1048
$ws.Cells.Item(35,4).Value = -216
$ws.Cells.Item(35,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(35,7).Value = "0.79 : 1"

$ws.Cells.Item(34,2).Value = 140
$ws.Cells.Item(34,3).Value = 109
$ws.Cells.Item(34,4).Value = 31
$ws.Cells.Item(34,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(34,7).Value = "1.28 : 1"

$ws.Cells.Item(33,2).Value = 273
$ws.Cells.Item(33,3).Value = 535
$ws.Cells.Item(33,4).Value = -262
$ws.Cells.Item(33,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(33,7).Value = "0.51 : 1"

$ws.Cells.Item(32,2).Value = 497
$ws.Cells.Item(32,3).Value = 541
$ws.Cells.Item(32,4).Value = -44
$ws.Cells.Item(32,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(32,7).Value = "0.92 : 1"

$ws.Cells.Item(31,2).Value = 75
$ws.Cells.Item(31,3).Value = 352
$ws.Cells.Item(31,4).Value = -277
$ws.Cells.Item(31,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(31,7).Value = "0.21 : 1"

$ws.Cells.Item(30,2).Value = 38
$ws.Cells.Item(30,3).Value = 6
$ws.Cells.Item(30,4).Value = 32
$ws.Cells.Item(30,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(30,7).Value = "6.33 : 1"

$ws.Cells.Item(29,2).Value = 583
$ws.Cells.Item(29,3).Value = 402
$ws.Cells.Item(29,4).Value = 181
$ws.Cells.Item(29,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(29,7).Value = "1.45 : 1"

$ws.Cells.Item(28,2).Value = 75
$ws.Cells.Item(28,3).Value = 107
$ws.Cells.Item(28,4).Value = -32
$ws.Cells.Item(28,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(28,7).Value = "0.70 : 1"

$ws.Cells.Item(27,2).Value = 201
$ws.Cells.Item(27,3).Value = 206
$ws.Cells.Item(27,4).Value = -5
$ws.Cells.Item(27,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(27,7).Value = "0.98 : 1"

$ws.Cells.Item(26,2).Value = 0
$ws.Cells.Item(26,3).Value = 0
$ws.Cells.Item(26,4).Value = 0

$ws.Cells.Item(25,2).Value = 136
$ws.Cells.Item(25,3).Value = 294
$ws.Cells.Item(25,4).Value = -158
$ws.Cells.Item(25,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(25,7).Value = "0.46 : 1"

$ws.Cells.Item(24,2).Value = 1889
$ws.Cells.Item(24,3).Value = 1154
$ws.Cells.Item(24,4).Value = 735
$ws.Cells.Item(24,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(24,7).Value = "1.64 : 1"

$ws.Cells.Item(23,2).Value = 644
$ws.Cells.Item(23,3).Value = 339
$ws.Cells.Item(23,4).Value = 305
$ws.Cells.Item(23,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(23,7).Value = "1.90 : 1"

$ws.Cells.Item(22,2).Value = 18
$ws.Cells.Item(22,3).Value = 147
$ws.Cells.Item(22,4).Value = -129
$ws.Cells.Item(22,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(22,7).Value = "0.12 : 1"

$ws.Cells.Item(21,2).Value = 414
$ws.Cells.Item(21,3).Value = 390
$ws.Cells.Item(21,4).Value = 24
$ws.Cells.Item(21,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(21,7).Value = "1.06 : 1"

$ws.Cells.Item(20,2).Value = 1
$ws.Cells.Item(20,3).Value = 119
$ws.Cells.Item(20,4).Value = -118
$ws.Cells.Item(20,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(20,7).Value = "0.01 : 1"

$ws.Cells.Item(19,2).Value = 503
$ws.Cells.Item(19,3).Value = 428
$ws.Cells.Item(19,4).Value = 75
$ws.Cells.Item(19,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(19,7).Value = "1.18 : 1"

$ws.Cells.Item(18,2).Value = 76
$ws.Cells.Item(18,3).Value = 86
$ws.Cells.Item(18,4).Value = -10
$ws.Cells.Item(18,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(18,7).Value = "0.88 : 1"

$ws.Cells.Item(17,2).Value = 581
$ws.Cells.Item(17,3).Value = 429
$ws.Cells.Item(17,4).Value = 152
$ws.Cells.Item(17,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(17,7).Value = "1.35 : 1"

$ws.Cells.Item(16,2).Value = 82
$ws.Cells.Item(16,3).Value = 144
$ws.Cells.Item(16,4).Value = -62
$ws.Cells.Item(16,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(16,7).Value = "0.57 : 1"

$ws.Cells.Item(15,2).Value = 49
$ws.Cells.Item(15,3).Value = 117
$ws.Cells.Item(15,4).Value = -68
$ws.Cells.Item(15,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(15,7).Value = "0.42 : 1"

$ws.Cells.Item(14,2).Value = 131
$ws.Cells.Item(14,3).Value = 272
$ws.Cells.Item(14,4).Value = -141
$ws.Cells.Item(14,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(14,7).Value = "0.48 : 1"

$ws.Cells.Item(13,2).Value = 205
$ws.Cells.Item(13,3).Value = 103
$ws.Cells.Item(13,4).Value = 102
$ws.Cells.Item(13,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(13,7).Value = "1.99 : 1"

$ws.Cells.Item(12,2).Value = 86
$ws.Cells.Item(12,3).Value = 15
$ws.Cells.Item(12,4).Value = 71
$ws.Cells.Item(12,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(12,7).Value = "5.73 : 1"

$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0

$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 57
$ws.Cells.Item(10,4).Value = -57
$ws.Cells.Item(10,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(10,7).Value = "0.00 : 1"

$ws.Cells.Item(9,2).Value = 56
$ws.Cells.Item(9,3).Value = 78
$ws.Cells.Item(9,4).Value = -22
$ws.Cells.Item(9,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(9,7).Value = "0.72 : 1"

$ws.Cells.Item(8,2).Value = 137
$ws.Cells.Item(8,3).Value = 200
$ws.Cells.Item(8,4).Value = -63
$ws.Cells.Item(8,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(8,7).Value = "0.69 : 1"

$ws.Cells.Item(7,2).Value = 231
$ws.Cells.Item(7,3).Value = 168
$ws.Cells.Item(7,4).Value = 63
$ws.Cells.Item(7,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(7,7).Value = "1.38 : 1"

$ws.Cells.Item(6,2).Value = 1039
$ws.Cells.Item(6,3).Value = 1587
$ws.Cells.Item(6,4).Value = -548
$ws.Cells.Item(6,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(6,7).Value = "0.65 : 1"

$ws.Cells.Item(5,2).Value = 55
$ws.Cells.Item(5,3).Value = 103
$ws.Cells.Item(5,4).Value = -48
$ws.Cells.Item(5,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(5,7).Value = "0.53 : 1"

$ws.Cells.Item(4,2).Value = 1233
$ws.Cells.Item(4,3).Value = 1240
$ws.Cells.Item(4,4).Value = -7
$ws.Cells.Item(4,6).Value = "We lent more than we borrowed"
$ws.Cells.Item(4,7).Value = "0.99 : 1"

$ws.Cells.Item(3,2).Value = 548
$ws.Cells.Item(3,3).Value = 498
$ws.Cells.Item(3,4).Value = 50
$ws.Cells.Item(3,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(3,7).Value = "1.10 : 1"

$ws.Cells.Item(2,2).Value = 1408
$ws.Cells.Item(2,3).Value = 1139
$ws.Cells.Item(2,4).Value = 269
$ws.Cells.Item(2,5).Value = "We borrowerd more than we lent"
$ws.Cells.Item(2,7).Value = "1.24 : 1"
